{"js": "// Office.js (Word JavaScript API) script.\n// This is the body of `async (context) => { ... }`.\n//\n// The edit:\n//  1) Three existing paragraphs (\"In an ideal situation...\", \"Such flexible\n//     architecture...\", \"Therefore, applications have a scaling limit...\")\n//     had their internal run/grammar-check (w:proofErr) structure cleaned up\n//     by Word on save; the visible text of those paragraphs is unchanged, so\n//     we simply re-assert their text (replace-with-self) which causes the\n//     runs to be rebuilt as a single clean run.\n//  2) A brand-new paragraph about fault tolerance is inserted right after\n//     the \"Therefore, applications have a scaling limit...\" paragraph (and\n//     before the trailing empty paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three paragraphs that need their run structure normalized, and\n// the paragraph after which the new content must be inserted, by matching\n// on a stable text prefix (robust to index drift).\nconst items = paragraphs.items;\n\nfunction findParagraph(startsWith) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(startsWith) === 0) {\n      return items[i];\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + startsWith);\n}\n\nconst pIdeal = findParagraph(\"In an ideal situation\");\nconst pSuch = findParagraph(\"Such flexible architecture\");\nconst pTherefore = findParagraph(\"Therefore, applications have a scaling limit\");\n\n// 1) Re-assert the text of each paragraph (no visible change) so the engine\n// rebuilds them as a single, clean run without leftover proofing markers.\npIdeal.insertText(pIdeal.text, Word.InsertLocation.replace);\nawait context.sync();\n\npSuch.insertText(pSuch.text, Word.InsertLocation.replace);\nawait context.sync();\n\npTherefore.insertText(pTherefore.text, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Insert the new \"fault tolerance\" paragraph right after pTherefore.\nconst newPara = pTherefore.insertParagraph(\n  \"In conclusion, let's talk about fault tolerance. This is a measure of how the system \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nnewPara.insertText(\"resists\", Word.InsertLocation.end);\nawait context.sync();\n\nnewPara.insertText(\n  \" failure. If we are talking about a system managed by a load balancer, a catastrophe will not happen. The failure of one instance will entail an increase in the load on the remaining ones. And if the autoscaling mode is enabled, the system itself will create another instance to replace the faulty one. The situation is much worse in the case of old systems, where there is only one server or a small number of them. Often, maintenance (restart) of such systems is in manual or semi-automatic mode. The failure of the server will entail the failure of the \",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\nnewPara.insertText(\"system\", Word.InsertLocation.end);\nawait context.sync();\n\nnewPara.insertText(\". Which often happens.\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# The edit:\n#  1) Three existing paragraphs (\"In an ideal situation...\", \"Such flexible\n#     architecture...\", \"Therefore, applications have a scaling limit...\")\n#     had their internal run/grammar-check (w:proofErr) structure cleaned up\n#     by Word on save; the visible text is unchanged, so we force Word to\n#     rebuild each paragraph's runs as a single clean run by appending a\n#     marker character and then removing it again.\n#  2) A brand-new paragraph about fault tolerance is inserted right after\n#     the \"Therefore, applications have a scaling limit...\" paragraph (and\n#     before the trailing empty paragraph).\n\n$doc = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($prefix) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $txt = $doc.Paragraphs($i).Range.Text\n        if ($txt.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Normalize-Paragraph($prefix) {\n    $idx = Find-ParagraphIndex $prefix\n    $p = $doc.Paragraphs($idx)\n    $r = $doc.Range($p.Range.Start, $p.Range.End - 1)\n    $original = $r.Text\n    # Append a throwaway marker, then strip it again. This forces the\n    # paragraph's runs (and any w:proofErr grammar markers) to be rebuilt\n    # from scratch as a single clean run, while leaving the visible text\n    # exactly as it was.\n    $r.Text = $original + \"@\"\n    $idx2 = Find-ParagraphIndex $prefix\n    $p2 = $doc.Paragraphs($idx2)\n    $r2 = $doc.Range($p2.Range.Start, $p2.Range.End - 1)\n    $r2.Text = $r2.Text.Substring(0, $r2.Text.Length - 1)\n}\n\nNormalize-Paragraph \"In an ideal situation\"\nNormalize-Paragraph \"Such flexible architecture\"\nNormalize-Paragraph \"Therefore, applications have a scaling limit\"\n\n# 2) Insert the new \"fault tolerance\" paragraph right after the \"Therefore,\n# applications have a scaling limit...\" paragraph.\n$thereforeIdx = Find-ParagraphIndex \"Therefore, applications have a scaling limit\"\n$thereforePara = $doc.Paragraphs($thereforeIdx)\n$thereforePara.Range.InsertParagraphAfter()\n\n$newIdx = $thereforeIdx + 1\n$newPara = $doc.Paragraphs($newIdx)\n$newText = \"In conclusion, let's talk about fault tolerance. This is a measure of how the system resists failure. If we are talking about a system managed by a load balancer, a catastrophe will not happen. The failure of one instance will entail an increase in the load on the remaining ones. And if the autoscaling mode is enabled, the system itself will create another instance to replace the faulty one. The situation is much worse in the case of old systems, where there is only one server or a small number of them. Often, maintenance (restart) of such systems is in manual or semi-automatic mode. The failure of the server will entail the failure of the system. Which often happens.\"\n$newPara.Range.Text = $newText\n"}
